$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellAddr, $text)
    $ws.Range("A1").Formula = '="' + $text + '"'
    $ws.Range("A1").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $ws.Range("A1").Clear()
}

Set-TextValue $ws 'D2' '25.709.88'
Set-TextValue $ws 'E2' '  -4.08%  '
Set-TextValue $ws 'D3' '1.815.22'
Set-TextValue $ws 'E3' '  -2.96%  '
Set-TextValue $ws 'D4' '1.002'
Set-TextValue $ws 'E4' '  +0.08%  '
Set-TextValue $ws 'D5' '278.37'
Set-TextValue $ws 'E5' '  -7.54%  '
Set-TextValue $ws 'E6' '  +0.04%  '
Set-TextValue $ws 'D7' '0.5078'
Set-TextValue $ws 'E7' '  -4.89%  '
Set-TextValue $ws 'D8' '0.3532'
Set-TextValue $ws 'E8' '  -5.65%  '
Set-TextValue $ws 'D9' '44.37'
Set-TextValue $ws 'E9' '  -2.41%  '
Set-TextValue $ws 'D10' '0.06685'
Set-TextValue $ws 'E10' '  -6.98%  '
Set-TextValue $ws 'D11' '19.88'
Set-TextValue $ws 'E11' '  -8.11%  '
Set-TextValue $ws 'D12' '0.8223'
Set-TextValue $ws 'E12' '  -7.48%  '
Set-TextValue $ws 'D13' '0.07863'
Set-TextValue $ws 'E13' '  -3.73%  '
Set-TextValue $ws 'D14' '1.820.36'
Set-TextValue $ws 'E14' '  -2.79%  '
Set-TextValue $ws 'D15' '5.071'
Set-TextValue $ws 'E15' '  -4.47%  '
Set-TextValue $ws 'D16' '87.61'
Set-TextValue $ws 'E16' '  -5.64%  '
Set-TextValue $ws 'E17' '  +0.15%  '
Set-TextValue $ws 'D18' '14.06'
Set-TextValue $ws 'E18' '  -5.29%  '
Set-TextValue $ws 'E19' '  +0.06%  '
Set-TextValue $ws 'D20' '0.000008023'
Set-TextValue $ws 'E20' '  -5.70%  '
Set-TextValue $ws 'D21' '25.750.74'
Set-TextValue $ws 'E21' '  -4.04%  '
Set-TextValue $ws 'D22' '4.740'
Set-TextValue $ws 'E22' '  -4.98%  '
Set-TextValue $ws 'D23' '9.991'
Set-TextValue $ws 'E23' '  -5.89%  '
Set-TextValue $ws 'D24' '6.105'
Set-TextValue $ws 'E24' '  -4.34%  '
Set-TextValue $ws 'D25' '2.241'
Set-TextValue $ws 'E25' '  -3.28%  '
Set-TextValue $ws 'D26' '142.37'
Set-TextValue $ws 'E26' '  -2.54%  '
Set-TextValue $ws 'D27' '1.667'
Set-TextValue $ws 'E27' '  -3.62%  '
Set-TextValue $ws 'D29' '109.18'
Set-TextValue $ws 'E29' '  -4.22%  '
Set-TextValue $ws 'D30' '4.330'
Set-TextValue $ws 'E30' '  -8.31%  '
Set-TextValue $ws 'D31' '4.224'
Set-TextValue $ws 'E31' '  -8.82%  '
Set-TextValue $ws 'D32' '0.08737'
Set-TextValue $ws 'E32' '  -4.57%  '
Set-TextValue $ws 'D33' '0.04872'
Set-TextValue $ws 'E33' '  -3.13%  '
Set-TextValue $ws 'D34' '0.7276'
Set-TextValue $ws 'E34' '  -9.51%  '
Set-TextValue $ws 'B35' 'HuobiToken'
Set-TextValue $ws 'C35' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D35' '2.885'
Set-TextValue $ws 'E35' '  -2.00%  '
Set-TextValue $ws 'B36' 'ARBITRUM'
Set-TextValue $ws 'C36' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws 'D36' '1.131'
Set-TextValue $ws 'E36' '  -3.75%  '
Set-TextValue $ws 'D37' '3.148'
Set-TextValue $ws 'E37' '  -1.50%  '
Set-TextValue $ws 'D38' '2.367'
Set-TextValue $ws 'E38' '  -12.58%  '
Set-TextValue $ws 'D39' '0.01848'
Set-TextValue $ws 'E39' '  -5.33%  '
Set-TextValue $ws 'D40' '0.5136'
Set-TextValue $ws 'E40' '  -16.22%  '
Set-TextValue $ws 'D41' '0.9697'
Set-TextValue $ws 'E41' '  -8.87%  '
Set-TextValue $ws 'D42' '114.13'
Set-TextValue $ws 'E42' '  -0.57%  '
Set-TextValue $ws 'D43' '6.224'
Set-TextValue $ws 'E43' '  -4.76%  '
Set-TextValue $ws 'D44' '7.998'
Set-TextValue $ws 'E44' '  -8.97%  '
Set-TextValue $ws 'E45' '  +0.08%  '
Set-TextValue $ws 'D46' '0.4525'
Set-TextValue $ws 'E46' '  -13.57%  '
Set-TextValue $ws 'D47' '0.1366'
Set-TextValue $ws 'E47' '  -8.44%  '
Set-TextValue $ws 'D48' '36.36'
Set-TextValue $ws 'E48' '  -3.34%  '
Set-TextValue $ws 'D49' '9.164'
Set-TextValue $ws 'E49' '  -7.95%  '
Set-TextValue $ws 'D50' '1.502'
Set-TextValue $ws 'E50' '  -8.92%  '
Set-TextValue $ws 'D51' '0.05824'
Set-TextValue $ws 'E51' '  -3.83%  '
